$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# K3: "2022-01-18" (text) -> real date serial with custom format "yyyy. M. d"
$ws.Range("K3").NumberFormat = "yyyy. M. d"
$ws.Range("K3").Value = 44579

# K4: "2021-10-25" (text) -> real date serial with custom format "yyyy. M. d"
$ws.Range("K4").NumberFormat = "yyyy. M. d"
$ws.Range("K4").Value = 44494

# K5: "2021-09-21" -> "2021-07-10" (stays text, matches the deduplicated K6 text)
$ws.Range("K5").Value = "2021-07-10"

# K9: "2021-08-02" -> "2021.08.02" (dot-separated, stays text)
$ws.Range("K9").Value = "2021.08.02"

# K10: "2022-01-21" -> "2022/01/21" (slash-separated, stays text)
$ws.Range("K10").Value = "2022/01/21"

Write-Host "done"
